$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was ID 3 / n4zdfr8rz4cdlxzf49uy / 08-07-2024 / invoices/12.jpg)
# Column A holds a plain-text row number ("1") that Excel would otherwise
# auto-detect as a number, so force Text format, write it, then clear the
# format again so the cell keeps its original (default) style but the
# stored value stays a literal string.
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "1"
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,2).Value = "bzfw5m88g0fymg6lk5cle"
$ws.Cells.Item(2,3).Value = "16-07-2024"
$ws.Cells.Item(2,4).Value = "https://rpachallengeocr.azurewebsites.net/invoices/12.jpg"

# Row 3 (was ID 9 / j6qky8ysjflms7kciqj97i / 14-07-2024 / invoices/8.jpg)
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2"
$ws.Cells.Item(3,1).ClearFormats()
$ws.Cells.Item(3,2).Value = "agkf70jn9satt1rxtvyy6"
# "07-07-2024" is itself a valid date (day/month both <=12), so Excel would
# silently convert it to a date serial number unless forced to Text too.
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value = "07-07-2024"
$ws.Cells.Item(3,3).ClearFormats()
$ws.Cells.Item(3,4).Value = "https://rpachallengeocr.azurewebsites.net/invoices/8.jpg"

# Row 4 (was ID 12 / g15db3dv9zupp579hzbzm / 15-07-2024 / invoices/7.jpg)
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "3"
$ws.Cells.Item(4,1).ClearFormats()
$ws.Cells.Item(4,2).Value = "wj34k48z92mgkik0lpdt1g"
$ws.Cells.Item(4,3).Value = "25-07-2024"
$ws.Cells.Item(4,4).Value = "https://rpachallengeocr.azurewebsites.net/invoices/10.jpg"
